# Update "想去人数" (F column) counts on both the "展览" and "全部类型" sheets
# F2: 1339 -> 1340
# F3: 75   -> 76
# F4: 3    -> 5

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1340
    $ws.Range("F3").Value = 76
    $ws.Range("F4").Value = 5
}
